# Add a new "Begründung" / "Absicherung von Skripten" row to the
# requirements table on Sheet1 (row 6), then move the active selection
# to B7, mirroring the state Excel leaves the sheet in after the user
# typed into B6 and pressed Enter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Begründung"
$ws.Range("B6").Value = "Absicherung von Skripten"

$ws.Range("B7").Select()
